$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.794.62"
$ws.Range("E2").Value = "  -4.37%  "
$ws.Range("D3").Value = "3.078.16"
$ws.Range("E3").Value = "  -3.52%  "
$ws.Range("E4").Value = "  -0.26%  "
$ws.Range("D5").Value = "'537.84"
$ws.Range("E5").Value = "  -5.44%  "
$ws.Range("D6").Value = "'132.83"
$ws.Range("E6").Value = "  -10.92%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "3.074.08"
$ws.Range("E8").Value = "  -3.37%  "
$ws.Range("D9").Value = "'0.488"
$ws.Range("E9").Value = "  -3.55%  "
$ws.Range("D10").Value = "'0.152"
$ws.Range("E10").Value = "  -4.71%  "
$ws.Range("D11").Value = "'6.16"
$ws.Range("E11").Value = "  -10.24%  "
$ws.Range("D12").Value = "'0.457"
$ws.Range("E12").Value = "  -4.76%  "
$ws.Range("D13").Value = "'0.0000224"
$ws.Range("E13").Value = "  -2.45%  "
$ws.Range("D14").Value = "'34.20"
$ws.Range("E14").Value = "  -9.74%  "
$ws.Range("D15").Value = "3.531.41"
$ws.Range("E15").Value = "  -4.75%  "
$ws.Range("D16").Value = "62.705.45"
$ws.Range("E16").Value = "  -4.77%  "
$ws.Range("E17").Value = "  -3.15%  "
$ws.Range("D18").Value = "3.069.09"
$ws.Range("E18").Value = "  -4.29%  "
$ws.Range("D19").Value = "'6.57"
$ws.Range("E19").Value = "  -6.70%  "
$ws.Range("D20").Value = "'478.80"
$ws.Range("E20").Value = "  -9.87%  "
$ws.Range("D21").Value = "'13.25"
$ws.Range("E21").Value = "  -8.04%  "
$ws.Range("D22").Value = "'0.696"
$ws.Range("E22").Value = "  -5.23%  "
$ws.Range("D23").Value = "'7.14"
$ws.Range("E23").Value = "  -6.55%  "
$ws.Range("D24").Value = "'78.51"
$ws.Range("E24").Value = "  -2.36%  "
$ws.Range("D25").Value = "'12.00"
$ws.Range("E25").Value = "  -9.78%  "
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "  +0.45%  "
$ws.Range("D27").Value = "'2.69"
$ws.Range("E27").Value = "  -7.18%  "
$ws.Range("D28").Value = "'8.13"
$ws.Range("E28").Value = "  -11.58%  "
$ws.Range("D29").Value = "'0.997"
$ws.Range("E29").Value = "  -0.53%  "
$ws.Range("D30").Value = "'25.75"
$ws.Range("E30").Value = "  -5.33%  "
$ws.Range("D31").Value = "'1.87"
$ws.Range("E31").Value = "  -16.05%  "
$ws.Range("D32").Value = "'1.10"
$ws.Range("E32").Value = "  -5.27%  "
$ws.Range("D33").Value = "'58.21"
$ws.Range("E33").Value = "  +7.04%  "
$ws.Range("D34").Value = "'2.39"
$ws.Range("E34").Value = "  -11.20%  "
$ws.Range("D35").Value = "'5.92"
$ws.Range("E35").Value = "  -5.46%  "
$ws.Range("D36").Value = "'5.19"
$ws.Range("E36").Value = "  -6.41%  "
$ws.Range("D37").Value = "'473.60"
$ws.Range("E37").Value = "  -13.84%  "
$ws.Range("D38").Value = "3.113.47"
$ws.Range("E38").Value = "  -2.29%  "
$ws.Range("D39").Value = "'0.0390"
$ws.Range("E39").Value = "  -12.28%  "
$ws.Range("D40").Value = "'0.0788"
$ws.Range("E40").Value = "  -6.81%  "
$ws.Range("B41").Value = "Cosmos"
$ws.Range("C41").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D41").Value = "'8.03"
$ws.Range("E41").Value = "  -5.47%  "
$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D42").Value = "'0.113"
$ws.Range("E42").Value = "  -10.97%  "
$ws.Range("D43").Value = "'2.60"
$ws.Range("E43").Value = "  -8.79%  "
$ws.Range("D44").Value = "'1.00"
$ws.Range("E44").Value = "  +0.03%  "
$ws.Range("D45").Value = "'0.249"
$ws.Range("E45").Value = "  -10.01%  "
$ws.Range("D46").Value = "'2.01"
$ws.Range("E46").Value = "  -12.14%  "
$ws.Range("D47").Value = "'24.32"
$ws.Range("E47").Value = "  -8.12%  "
$ws.Range("D48").Value = "'117.96"
$ws.Range("E48").Value = "  -4.44%  "
$ws.Range("D49").Value = "'0.107"
$ws.Range("E49").Value = "  -3.46%  "
$ws.Range("D50").Value = "0.0₃0508"
$ws.Range("E50").Value = "  -6.09%  "
$ws.Range("D51").Value = "'1.98"
$ws.Range("E51").Value = "  -8.81%  "
